# Unbind PIP component sheet: add an English (column C) translation column
# next to the existing key (A) / Chinese text (B) columns, and mark the
# "IP" substrings inside the Chinese text with the Arial font (matching the
# mixed CJK/Latin run formatting used throughout the translated workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: keys (unchanged, just re-assert them) ------------------
$ws.Range("A1").Value = "unbind_tips_one"
$ws.Range("A2").Value = "unbind_tips_two"
$ws.Range("A3").Value = "confirm_unbind_pip"

# --- Column B: Chinese text (unchanged content) ------------------------
$ws.Range("B1").Value = "解绑后若需绑定公网IP，仅支持绑定与云主机同可用区或带全可用区属性的公网IP"
$ws.Range("B2").Value = "解绑后若需绑定公网IP，仅支持绑定全可用区属性的公网IP"
$ws.Range("B3").Value = "确认解绑公网IP"

# --- Column C (new): English translation --------------------------------
$ws.Range("C1").Value = "If the EIP needs to be associated after disassociation, it only supports associating the EIP that is in the same availability zone with the virtual machine or that has the full availability zone attribute"
$ws.Range("C2").Value = "If the EIP needs to be associated after disassociation, it only supports associating EIP that has the full availability zone attribute"
$ws.Range("C3").Value = "Confirm to disassociate the EIP"

# --- Rich-text runs in column B: "IP" segments use Arial, the rest 宋体 --
$ws.Range("B1").Characters(1, 9).Font.Name = "宋体"
$ws.Range("B1").Characters(1, 9).Font.Size = 12
$ws.Range("B1").Characters(10, 2).Font.Name = "Arial"
$ws.Range("B1").Characters(10, 2).Font.Size = 12
$ws.Range("B1").Characters(12, 25).Font.Name = "宋体"
$ws.Range("B1").Characters(12, 25).Font.Size = 12
$ws.Range("B1").Characters(37, 2).Font.Name = "Arial"
$ws.Range("B1").Characters(37, 2).Font.Size = 12

$ws.Range("B2").Characters(1, 9).Font.Name = "宋体"
$ws.Range("B2").Characters(1, 9).Font.Size = 12
$ws.Range("B2").Characters(10, 2).Font.Name = "Arial"
$ws.Range("B2").Characters(10, 2).Font.Size = 12
$ws.Range("B2").Characters(12, 15).Font.Name = "宋体"
$ws.Range("B2").Characters(12, 15).Font.Size = 12
$ws.Range("B2").Characters(27, 2).Font.Name = "Arial"
$ws.Range("B2").Characters(27, 2).Font.Size = 12

$ws.Range("B3").Characters(1, 6).Font.Name = "宋体"
$ws.Range("B3").Characters(1, 6).Font.Size = 12
$ws.Range("B3").Characters(7, 2).Font.Name = "Arial"
$ws.Range("B3").Characters(7, 2).Font.Size = 12

# --- Whole used range: Arial 12pt as the cell-level (non-rich) font -----
$ws.Range("A1:C3").Font.Name = "Arial"
$ws.Range("A1:C3").Font.Size = 12

# --- Column widths (characters); engine applies a fixed +5/6 pixel pad -
$ws.Columns.Item(1).ColumnWidth = 19.666666666666668   # -> stored 20.5
$ws.Columns.Item(2).ColumnWidth = 79.66666666666667    # -> stored 80.5
$ws.Columns.Item(3).ColumnWidth = 226.04166666666666   # -> stored ~226.875

# --- Selection / active cell, matches the saved view state -------------
$ws.Range("B12").Select()
